$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regular cell value updates ---
$ws.Cells.Item(2, 4).Value = '38.265.23'
$ws.Cells.Item(2, 5).Value = '  +3.12%  '

$ws.Cells.Item(3, 4).Value = '2.067.37'
$ws.Cells.Item(3, 5).Value = '  +2.38%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '230.43'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.21%  '

$ws.Cells.Item(6, 5).Value = '  +0.97%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '61.27'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +9.79%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.14%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.387'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +3.30%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0807'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +3.70%  '

$ws.Cells.Item(11, 5).Value = '  +1.85%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '14.88'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +4.97%  '

$ws.Cells.Item(13, 4).Value = '2.372.80'
$ws.Cells.Item(13, 5).Value = '  +2.35%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '21.32'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +6.87%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.766'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +3.67%  '

$ws.Cells.Item(16, 5).Value = '  +2.11%  '

$ws.Cells.Item(17, 4).Value = '2.068.12'
$ws.Cells.Item(17, 5).Value = '  +2.46%  '

$ws.Cells.Item(18, 4).Value = '38.181.16'
$ws.Cells.Item(18, 5).Value = '  +3.00%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.28'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.59%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '70.19'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.70%  '

$ws.Cells.Item(21, 5).Value = '  +2.93%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '226.40'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +1.20%  '

$ws.Cells.Item(23, 5).Value = '  -0.13%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.42'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.76%  '

$ws.Cells.Item(25, 5).Value = '  +2.11%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '166.30'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.00%  '

$ws.Cells.Item(27, 5).Value = '  +2.83%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.135'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +6.45%  '

$ws.Cells.Item(29, 5).Value = '  +1.70%  '

$ws.Cells.Item(30, 5).Value = '  +0.15%  '

$ws.Cells.Item(31, 5).Value = '  +1.82%  '

$ws.Cells.Item(32, 5).Value = '  +3.13%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.62'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +3.51%  '

$ws.Cells.Item(34, 5).Value = '  +9.63%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0606'
$ws.Cells.Item(35, 4).Style = "Normal"

$ws.Cells.Item(38, 5).Value = '  +4.99%  '

$ws.Cells.Item(39, 5).Value = '  -0.01%  '

$ws.Cells.Item(40, 4).Value = '1.523.38'
$ws.Cells.Item(40, 5).Value = '  +3.77%  '

$ws.Cells.Item(43, 5).Value = '  +2.22%  '

$ws.Cells.Item(44, 5).Value = '  +4.25%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0928'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.39%  '

$ws.Cells.Item(46, 5).Value = '  +1.74%  '

$ws.Cells.Item(47, 5).Value = '  -2.22%  '

$ws.Cells.Item(49, 5).Value = '  +1.18%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.11'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.26%  '

$ws.Cells.Item(51, 4).Value = '2.262.15'
$ws.Cells.Item(51, 5).Value = '  +2.60%  '

# --- Swapped rows (coin order changed) ---
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.33'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.88%  '

$ws.Cells.Item(37, 2).Value = 'THORChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.26'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +13.94%  '

$ws.Cells.Item(41, 2).Value = 'Aave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '98.34'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +3.94%  '

$ws.Cells.Item(42, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '17.11'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +6.24%  '
